$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers ---
$ws.Range("H1").Value = "HOSTNAME"
$ws.Range("I1").Value = "EXPL_QUERY1"
$ws.Range("J1").Value = "EXPL_QUERY2"
$ws.Range("K1").Value = "EXPL_QUERY3"
$ws.Range("L1").Value = "QUERY1"
$ws.Range("M1").Value = "QUERY2"
$ws.Range("N1").Value = "QUERY3"
$ws.Range("O1").Value = "DATABASE_1"
$ws.Range("P1").Value = "DATABASE_2"

# --- Row 2 plain values (O2/P2 first, then H2, to match shared-string insertion order) ---
$ws.Range("O2").Value = "Digisales_KPI"
$ws.Range("P2").Value = "Digisales_SAPM"
$ws.Range("H2").Value = "192.168.232.6"

# --- Row 2 formulas ---
$ws.Range("I2").Formula = '="Menampilkan Semua Data pada dbo.TMP_BOOSTER di Database " & O2'
$ws.Range("J2").Formula = '="Menampilkan Semua Data pada dbo.TMP_BOOSTER di Database " & P2'
$ws.Range("L2").Formula = '="USE " & O2 & "; Select * From dbo.TMP_BOOSTER"'
$ws.Range("M2").Formula = '="USE " & P2 & "; Select * From dbo.TMP_BOOSTER"'

# --- Row 2 formatting ---
# H2/O2/P2 -> same look as A2 (default font, vertical center, no wrap)
$ws.Range("A2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("O2").PasteSpecial(-4122)
$ws.Range("P2").PasteSpecial(-4122)

# I2/J2 -> same look as C2/D2/E2 (default font, vertical center + wrap)
$ws.Range("C2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("J2").PasteSpecial(-4122)

# K2/L2/M2 -> the Arial-10-black font + vertical center style already used by N2 (before edit)
$ws.Range("N2").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("M2").PasteSpecial(-4122)
$ws.Range("L2").WrapText = $true
$ws.Range("M2").WrapText = $true

# N2 no longer exists in the new layout
$ws.Range("N2").Clear()

# --- Sheet view changes ---
$ws.Range("I2").Select()
$ws.ActiveWindow.ScrollColumn = 5
